# "faster Mongolia cold stations min max retrieve"
#
# Renames the second sheet (plain "Sheet1") to "各省最低温" (per-province
# record-low temperatures) and switches the active/visible tab from the
# first sheet ("国家站极端温度") to this one, scrolled/selected near the
# bottom of its data (around row 54-61) instead of its previous J9 view.

$wb = $excel.ActiveWorkbook

$sheetExtremes = $wb.Worksheets.Item("国家站极端温度")
$sheetLows     = $wb.Worksheets.Item("Sheet1")

# Rename "Sheet1" -> "各省最低温"
$sheetLows.Name = "各省最低温"

# Leave sheet 1's own selection/scroll state untouched - it keeps showing
# L1439 - and just make the renamed sheet the active tab in the workbook.
$sheetLows.Activate()

# Scroll so row 54 is at the top of the viewport and select H61, matching
# the new saved view position on the renamed sheet.
$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
$sheetLows.Range("H61").Select() | Out-Null
